$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): rename / reorder columns, add new "Codigo" + "Responsavel"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Codigo"
$ws.Range("B1").Value = "Razao Social"
$ws.Range("C1").Value = "CNPJ"
$ws.Range("D1").Value = "Nome Fantasia"
$ws.Range("E1").Value = "E-mail"
$ws.Range("F1").Value = "Telefone"
$ws.Range("G1").Value = "Endereco"
$ws.Range("H1").Value = "Responsavel"

# ---------------------------------------------------------------------------
# Data row 2 - replace sample data (keep CNPJ digits-only, stored as text)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "EMP001"
$ws.Range("B2").Value = "Empresa Exemplo Ltda"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "12345678000190"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "Exemplo"
$ws.Range("E2").Value = "contato@exemplo.com"
$ws.Range("F2").Value = "(11) 99999-9999"
$ws.Range("G2").Value = "Rua Exemplo, 123 - Sao Paulo - SP"
$ws.Range("H2").Value = "Joao Silva"

# Remove the old second data row ("Outra Empresa S.A.")
$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Column widths (8 columns now)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(2).ColumnWidth = 39.166666666666664
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668
$ws.Columns.Item(7).ColumnWidth = 39.166666666666664
$ws.Columns.Item(8).ColumnWidth = 24.166666666666668

# ---------------------------------------------------------------------------
# Header style: bold white text on dark-blue fill, centered + wrapped
# (366092 RGB -> OLE/BGR color int 9592886)
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:H1")
$headerRange.Interior.Color = 9592886
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true

# ---------------------------------------------------------------------------
# Freeze header row (split below row 1) - select A2 to freeze, then back to A1
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# AutoFilter across the header row + matching hidden _FilterDatabase name
# ---------------------------------------------------------------------------
$ws.Range("A1:H1").AutoFilter()

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Empresas'!`$A`$1:`$H`$1")
$filterName.Visible = $false
